$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.171683
$ws.Range("H2").Value = 0.515049
$ws.Range("I2").Value = 0.05260263278194677
$ws.Range("J2").Value = 0.05260263278194677
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.110028
$ws.Range("N2").Value = 0.330084
$ws.Range("Q2").Value = 0.018889937124
$ws.Range("R2").Value = 0.170009434116
$ws.Range("S2").Value = 0.05260263278194677
$ws.Range("T2").Value = 0.05260263278194677

# Row 3
$ws.Range("I3").Value = 0.2351975975445997
$ws.Range("J3").Value = 0.2351975975445997
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.110028
$ws.Range("N3").Value = 0.330084
$ws.Range("Q3").Value = 0.08446094034400001
$ws.Range("R3").Value = 0.760148463096
$ws.Range("S3").Value = 0.2351975975445997
$ws.Range("T3").Value = 0.2351975975445997

# Row 4
$ws.Range("G4").Value = 1.354395
$ws.Range("H4").Value = 4.063185
$ws.Range("I4").Value = 0.4149784359936907
$ws.Range("J4").Value = 0.4149784359936907
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110028
$ws.Range("N4").Value = 0.330084
$ws.Range("Q4").Value = 0.14902137306
$ws.Range("R4").Value = 1.34119235754
$ws.Range("S4").Value = 0.4149784359936907
$ws.Range("T4").Value = 0.4149784359936907

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.5285906666666667
$ws.Range("H5").Value = 1.585772
$ws.Range("I5").Value = 0.1619569831062545
$ws.Range("J5").Value = 0.1619569831062546
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.110028
$ws.Range("N5").Value = 0.330084
$ws.Range("Q5").Value = 0.058159773872
$ws.Range("R5").Value = 0.523437964848
$ws.Range("S5").Value = 0.1619569831062545
$ws.Range("T5").Value = 0.1619569831062546

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.441472
$ws.Range("H6").Value = 1.324416
$ws.Range("I6").Value = 0.1352643505735082
$ws.Range("J6").Value = 0.1352643505735082
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.110028
$ws.Range("N6").Value = 0.330084
$ws.Range("Q6").Value = 0.04857428121600001
$ws.Range("R6").Value = 0.437168530944
$ws.Range("S6").Value = 0.1352643505735082
$ws.Range("T6").Value = 0.1352643505735082
